$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d18 = @'
Objective
id="Par1">The aim of this study was to investigate the clinical characteristics of Corona Virus Disease 2019 in Taizhou, China.


Methods
id="Par2">A single center retrospective observational study was performed between Jan 1, 2020 and Mar 11, 2020 at Taizhou Public Health Medical Center, Zhejiang, China.

 All patients with confirmed Corona Virus Disease 2019 were enrolled, and their clinical data were gathered by reviewing electronic medical records.

 Outcomes of severely ill patients and non-severely ill patients were compared.


Results
id="Par3">Of 145 hospitalized patients with COVID-19, the average age was 47.5 years old (standard deviation, 14.6) and 54.5% were men.

 Hypertension was the most common comorbidity (15.2%), followed by diabetes mellitus (9.7%).

 Common symptoms included dry cough (81.4%), fever (75.2%), anorexia (42.8%), fatigue (40.7%), chest tightness (32.4%), diarrhea (26.9%) and dizziness (20%).

 According to imaging examination, 79.3% patients showed bilateral pneumonia, 18.6% showed unilateral pneumonia, 61.4% showed ground-glass opacity, and 2.1% showed no abnormal result.

 Compared with non-severely ill patients, severely ill patients were older (mean, years, 52.8 vs.

 45.3, p &lt; 0.01), had a higher proportion of diabetes mellitus (16.3% vs.

 6.9%, p = 0.08), had a higher body mass index (mean, 24.78 vs.

 23.20, p = 0.02) and were more likely to have fever (90.7% vs.

 68.6%, p = 0.01), anorexia (60.5% vs.

 35.3%, p = 0.01), chest tightness (60.5% vs.

20.6%, p &lt; 0.01) and dyspnea (7.0% vs.

 0%, p = 0.03).

 Of the 43 severely ill patients, 6 (14%) received high-flow nasal cannula oxygen therapy, and 1 (2.3%) received invasive mechanical ventilation.


Conclusions
id="Par4">Older patients or patients with comorbidities such as obesity or diabetes mellitus were more likely to have severe condition.

 Treatments of COVID-19 is still experimental and more clinical trials are needed.



'@

$e18 = @'
[Qingqing%Chen%NULL%0,    Zhencang%Zheng%NULL%1,    Chao%Zhang%NULL%1,    Xijiang%Zhang%NULL%1,    Huijuan%Wu%NULL%1,    Jingdong%Wang%NULL%1,    Shuwei%Wang%NULL%1,    Cheng%Zheng%dr.zhengcheng@foxmail.com%1]
'@

$ws.Range("D18").Value = $d18
$ws.Range("E18").Value = $e18
$ws.Rows.Item(18).EntireRow.AutoFit()
